$wb = $excel.ActiveWorkbook

# Rename "Лист2" -> "GirHun"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "GirHun"

# Fill in the new table on GirHun (order chosen so new shared strings are
# appended in the same sequence as the target file: github URL, Type, Hint)
$ws2.Range("A1").Value = "URL"
$ws2.Range("A2").Value = "https://github.com/hqztrue/LeetCodeSolutions/"
$ws2.Range("B1").Value = "Type"
$ws2.Range("B2").Value = "Hint"

# Auto-fit column A to the new content
$ws2.Columns.Item(1).AutoFit() | Out-Null

# Make GirHun the active/selected sheet, with B2 selected
$ws2.Activate()
$ws2.Range("B2").Select()
